$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.190.03"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "3.614.81"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "194.29"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.57%  "
$ws.Range("D7").Value = "3.611.13"
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.621"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.679"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.62"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  +10.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.05"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "4.191.32"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("D16").Value = "3.618.56"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.53"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Value = "68.102.08"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.54"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "403.66"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +23.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +9.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.95"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.61"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.14"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.12"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +15.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.18"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.69"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "690.23"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +14.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.24"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.82"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.67"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.418"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.17%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "0.0₃0791"
$ws.Range("E40").Value = "  +5.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +16.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.11"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.33%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "3.145.28"
$ws.Range("E44").Value = "  +14.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0424"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.84"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.12"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.57"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("E51").Value = "  +1.51%  "
